# Generate Report for Handback
# Updates the handoff/handback timestamps (and the Overview "Latest HO Xliff
# Generate Date") for the second tracked file
# (8ac2b10a-90e9-4c00-a003-c701f467e106.md) now that it has been handed back.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is 8ac2b10a-90e9-4c00-a003-c701f467e106.md, column G is
# "Latest HO Xliff Generate Date"
$overview.Range("G3").Value = "2016-08-29 02:47:46"

# zh-cn sheet: row 3 is 8ac2b10a-90e9-4c00-a003-c701f467e106.md
#   column H = Correspond Handoff Datetime
#   column K = Correspond Handback DateTime
$zhcn.Range("H3").Value = "2016-08-29 02:47:41"
$zhcn.Range("K3").Value = "2016-08-29 02:48:15"

# de-de sheet: row 3 is 8ac2b10a-90e9-4c00-a003-c701f467e106.md
#   column H = Correspond Handoff Datetime
#   column K = Correspond Handback DateTime
$dede.Range("H3").Value = "2016-08-29 02:47:46"
$dede.Range("K3").Value = "2016-08-29 02:48:23"
